$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 12 (2020年) values: the original data had its decimal point
#     shifted 8 places to the left by mistake (e.g. 102095839 -> 1.02095839).
#     Restore the correct integer values. ---
$ws.Range("B12").Value = 102095839
$ws.Range("C12").Value = 425604368
$ws.Range("D12").Value = 14569041
$ws.Range("E12").Value = 342792219
$ws.Range("F12").Value = 42425631
$ws.Range("H12").Value = 480690099
$ws.Range("I12").Value = 7574980
$ws.Range("J12").Value = 2480282124
$ws.Range("K12").Value = 485012261
$ws.Range("L12").Value = 3151911008
$ws.Range("M12").Value = 561958065
$ws.Range("N12").Value = 87526798
$ws.Range("O12").Value = 30639227
$ws.Range("P12").Value = 14181667
$ws.Range("R12").Value = 76867754
$ws.Range("S12").Value = 1624639012
$ws.Range("T12").Value = 4322162

# --- Append new row 13 with 2021年 data ---
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 103937568
$ws.Range("C13").Value = 470520014
$ws.Range("D13").Value = 13973004
$ws.Range("E13").Value = 357780502
$ws.Range("F13").Value = 45075872
$ws.Range("H13").Value = 483034115
$ws.Range("I13").Value = 6747453
$ws.Range("J13").Value = 2574121355
$ws.Range("K13").Value = 487242918
$ws.Range("L13").Value = 3252812722
$ws.Range("M13").Value = 568006346
$ws.Range("N13").Value = 89964564
$ws.Range("O13").Value = 35697495
$ws.Range("P13").Value = 14115371
$ws.Range("R13").Value = 80540558
$ws.Range("S13").Value = 1650932101
$ws.Range("T13").Value = 4208803

# Copy the formatting of the 2020年 label cell (bordered, bold, centered
# header style) onto the new 2021年 label cell so row 13 matches the look
# of every other year row.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = "2021年"
